$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '28.930.10'
Set-TextValue 'E2' '  -0.09%  '
Set-TextValue 'D3' '1.920.01'
Set-TextValue 'E3' '  +0.91%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '324.45'
Set-TextValue 'E5' '  +0.06%  '
Set-TextValue 'E6' '  -0.02%  '
Set-TextValue 'D7' '0.4567'
Set-TextValue 'E7' '  -0.58%  '
Set-TextValue 'D8' '0.3801'
Set-TextValue 'E8' '  -0.28%  '
Set-TextValue 'D9' '0.07741'
Set-TextValue 'E9' '  +0.33%  '
Set-TextValue 'D10' '0.9755'
Set-TextValue 'E10' '  -0.54%  '
Set-TextValue 'D11' '22.30'
Set-TextValue 'E11' '  +1.18%  '
Set-TextValue 'D12' '1.910.17'
Set-TextValue 'E12' '  -3.09%  '
Set-TextValue 'D13' '5.695'
Set-TextValue 'E13' '  +0.47%  '
Set-TextValue 'D14' '6.952'
Set-TextValue 'E14' '  -0.20%  '
Set-TextValue 'D15' '0.06985'
Set-TextValue 'E15' '  -0.93%  '
Set-TextValue 'D16' '1.006'
Set-TextValue 'E16' '  +0.10%  '
Set-TextValue 'D17' '84.44'
Set-TextValue 'E17' '  +0.59%  '
Set-TextValue 'D18' '0.000009467'
Set-TextValue 'E18' '  -0.70%  '
Set-TextValue 'D19' '16.67'
Set-TextValue 'E19' '  -0.45%  '
Set-TextValue 'D20' '1.005'
Set-TextValue 'E20' '  +0.11%  '
Set-TextValue 'D21' '28.950.02'
Set-TextValue 'E21' '  +0.05%  '
Set-TextValue 'D22' '5.335'
Set-TextValue 'E22' '  +0.14%  '
Set-TextValue 'D23' '11.06'
Set-TextValue 'E23' '  +1.55%  '
Set-TextValue 'D24' '2.139.83'
Set-TextValue 'E24' '  -3.56%  '
Set-TextValue 'D25' '2.058'
Set-TextValue 'E25' '  -2.05%  '
Set-TextValue 'D26' '157.78'
Set-TextValue 'E26' '  +0.19%  '
Set-TextValue 'D27' '19.04'
Set-TextValue 'E27' '  -0.54%  '
Set-TextValue 'D28' '5.602'
Set-TextValue 'E28' '  +0.51%  '
Set-TextValue 'D29' '117.78'
Set-TextValue 'E29' '  +0.27%  '
Set-TextValue 'D30' '1.837'
Set-TextValue 'E30' '  -0.70%  '
Set-TextValue 'D31' '0.09269'
Set-TextValue 'E31' '  -0.17%  '
Set-TextValue 'D32' '0.8635'
Set-TextValue 'E32' '  +0.14%  '
Set-TextValue 'D33' '5.099'
Set-TextValue 'E33' '  +0.15%  '
Set-TextValue 'D34' '1.238'
Set-TextValue 'E34' '  -1.02%  '
Set-TextValue 'D35' '3.011'
Set-TextValue 'E35' '  -0.19%  '
Set-TextValue 'D36' '0.05683'
Set-TextValue 'E36' '  -0.15%  '
Set-TextValue 'D37' '1.150'
Set-TextValue 'E37' '  +0.34%  '
Set-TextValue 'D38' '1.003'
Set-TextValue 'E38' '  +0.03%  '
Set-TextValue 'D39' '0.02044'
Set-TextValue 'E39' '  +0.37%  '
Set-TextValue 'D40' '3.087'
Set-TextValue 'E40' '  +12.98%  '
Set-TextValue 'D41' '7.453'
Set-TextValue 'E41' '  -0.43%  '
Set-TextValue 'D42' '0.5487'
Set-TextValue 'E42' '  -0.41%  '
Set-TextValue 'D43' '0.1752'
Set-TextValue 'E43' '  -0.06%  '
Set-TextValue 'D44' '9.309'
Set-TextValue 'E44' '  +0.12%  '
Set-TextValue 'D45' '0.000002777'
Set-TextValue 'E45' '  +14.10%  '
Set-TextValue 'D46' '2.160'
Set-TextValue 'E46' '  +3.74%  '
Set-TextValue 'D47' '0.5164'
Set-TextValue 'E47' '  -0.60%  '
Set-TextValue 'D48' '0.06938'
Set-TextValue 'E48' '  +1.73%  '
Set-TextValue 'D49' '11.23'
Set-TextValue 'E49' '  -0.21%  '
Set-TextValue 'D50' '110.56'
Set-TextValue 'E50' '  -0.54%  '
Set-TextValue 'D51' '1.759'
Set-TextValue 'E51' '  -0.82%  '
